$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(28, 8).Value2 = 2196.8
$ws.Cells.Item(28, 9).Value2 = 970.2143
$ws.Cells.Item(28, 11).Value2 = 970.2143
$ws.Cells.Item(28, 13).Value2 = -485.2143

$ws.Cells.Item(55, 8).Value2 = 975.25
$ws.Cells.Item(55, 9).Value2 = 838.6667
$ws.Cells.Item(55, 10).Value2 = 1020.7778
$ws.Cells.Item(55, 11).Value2 = 838.6667
$ws.Cells.Item(55, 12).Value2 = 1020.7778
$ws.Cells.Item(55, 13).Value2 = -624.6667
$ws.Cells.Item(55, 14).Value2 = -1448.7778

$ws.Cells.Item(62, 8).Value2 = 6065.2856
$ws.Cells.Item(62, 9).Value2 = 2114.25
$ws.Cells.Item(62, 11).Value2 = 2114.25
$ws.Cells.Item(62, 13).Value2 = -1490.25

$ws.Cells.Item(65, 8).Value2 = 6065.2856
$ws.Cells.Item(65, 9).Value2 = 2114.25
$ws.Cells.Item(65, 11).Value2 = 10571.25
$ws.Cells.Item(65, 13).Value2 = -7451.25

$ws.Cells.Item(97, 8).Value2 = 556.75
$ws.Cells.Item(97, 10).Value2 = 556.75
$ws.Cells.Item(97, 12).Value2 = 1670.25
$ws.Cells.Item(97, 14).Value2 = -2662.25

$ws.Cells.Item(98, 8).Value2 = 403.7857
$ws.Cells.Item(98, 9).Value2 = 342.6154
$ws.Cells.Item(98, 11).Value2 = 342.6154
$ws.Cells.Item(98, 13).Value2 = 1155.3846

$ws.Cells.Item(100, 8).Value2 = 1464.1666
$ws.Cells.Item(100, 9).Value2 = 1457
$ws.Cells.Item(100, 10).Value2 = 1500
$ws.Cells.Item(100, 11).Value2 = 1457
$ws.Cells.Item(100, 12).Value2 = 1500
$ws.Cells.Item(100, 13).Value2 = -916
$ws.Cells.Item(100, 14).Value2 = -2582

$ws.Cells.Item(101, 8).Value2 = 474
$ws.Cells.Item(101, 9).Value2 = 335
$ws.Cells.Item(101, 11).Value2 = 1005
$ws.Cells.Item(101, 13).Value2 = 617

$ws.Cells.Item(103, 8).Value2 = 5492.857
$ws.Cells.Item(103, 9).Value2 = 7216.6665
$ws.Cells.Item(103, 11).Value2 = 21649.9995
$ws.Cells.Item(103, 13).Value2 = -21063.9995

$ws.Cells.Item(122, 8).Value2 = 403.7857
$ws.Cells.Item(122, 9).Value2 = 342.6154
$ws.Cells.Item(122, 11).Value2 = 1027.8462
$ws.Cells.Item(122, 13).Value2 = 1422.1538

$ws.Cells.Item(132, 8).Value2 = 10187.792
$ws.Cells.Item(132, 9).Value2 = 10295.772
$ws.Cells.Item(132, 11).Value2 = 30887.316
$ws.Cells.Item(132, 13).Value2 = -28357.316

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value2 = 1254.5454
$ws.Cells.Item(2, 9).Value2 = 1180
$ws.Cells.Item(2, 11).Value2 = 1180
$ws.Cells.Item(2, 13).Value2 = -1067

$ws.Cells.Item(61, 8).Value2 = 1502.6666
$ws.Cells.Item(61, 9).Value2 = 1502.6666
$ws.Cells.Item(61, 11).Value2 = 1502.6666
$ws.Cells.Item(61, 13).Value2 = -1290.6666

$ws.Cells.Item(74, 8).Value2 = 0
$ws.Cells.Item(74, 9).Value2 = 0
$ws.Cells.Item(74, 11).Value2 = 0
$ws.Cells.Item(74, 13).ClearContents()

$ws.Cells.Item(77, 8).Value2 = 0
$ws.Cells.Item(77, 9).Value2 = 0
$ws.Cells.Item(77, 11).Value2 = 0
$ws.Cells.Item(77, 13).ClearContents()

$ws.Cells.Item(97, 8).Value2 = 999.25
$ws.Cells.Item(97, 9).Value2 = 999.3333
$ws.Cells.Item(97, 11).Value2 = 999.3333
$ws.Cells.Item(97, 13).Value2 = -503.3333

$ws.Cells.Item(116, 8).Value2 = 1254.5454
$ws.Cells.Item(116, 9).Value2 = 1180
$ws.Cells.Item(116, 11).Value2 = 1180
$ws.Cells.Item(116, 13).Value2 = 1114

$ws.Cells.Item(122, 8).Value2 = 1333
$ws.Cells.Item(122, 9).Value2 = 1333
$ws.Cells.Item(122, 11).Value2 = 3999
$ws.Cells.Item(122, 13).Value2 = -1549

$ws.Cells.Item(132, 8).Value2 = 5000
$ws.Cells.Item(132, 9).Value2 = 5000
$ws.Cells.Item(132, 11).Value2 = 15000
$ws.Cells.Item(132, 13).Value2 = -12470

$ws.Cells.Item(136, 8).Value2 = 1502.6666
$ws.Cells.Item(136, 9).Value2 = 1502.6666
$ws.Cells.Item(136, 11).Value2 = 4507.9998
$ws.Cells.Item(136, 13).Value2 = -1957.9998

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value2 = 1254.5454
$ws.Cells.Item(3, 9).Value2 = 1180
$ws.Cells.Item(3, 11).Value2 = 1180
$ws.Cells.Item(3, 13).Value2 = -1066

$ws.Cells.Item(20, 8).Value2 = 1050
$ws.Cells.Item(20, 9).Value2 = 1050
$ws.Cells.Item(20, 11).Value2 = 1050
$ws.Cells.Item(20, 13).Value2 = -803

$ws.Cells.Item(86, 8).Value2 = 4889.8
$ws.Cells.Item(86, 9).Value2 = 1566.6666
$ws.Cells.Item(86, 10).Value2 = 6314
$ws.Cells.Item(86, 11).Value2 = 1566.6666
$ws.Cells.Item(86, 12).Value2 = 6314
$ws.Cells.Item(86, 13).Value2 = -443.6666
$ws.Cells.Item(86, 14).Value2 = -8560

$ws.Cells.Item(89, 8).Value2 = 4889.8
$ws.Cells.Item(89, 9).Value2 = 1566.6666
$ws.Cells.Item(89, 10).Value2 = 6314
$ws.Cells.Item(89, 11).Value2 = 7833.333000000001
$ws.Cells.Item(89, 12).Value2 = 31570
$ws.Cells.Item(89, 13).Value2 = -2217.333000000001
$ws.Cells.Item(89, 14).Value2 = -42802

$ws.Cells.Item(94, 8).Value2 = 730.1667
$ws.Cells.Item(94, 9).Value2 = 730.1667
$ws.Cells.Item(94, 11).Value2 = 730.1667
$ws.Cells.Item(94, 13).Value2 = -279.1667

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value2 = 952
$ws.Cells.Item(16, 9).Value2 = 502.5
$ws.Cells.Item(16, 11).Value2 = 502.5
$ws.Cells.Item(16, 13).Value2 = -215.5

$ws.Cells.Item(31, 8).Value2 = 6032.4287
$ws.Cells.Item(31, 9).Value2 = 1726.7
$ws.Cells.Item(31, 10).Value2 = 9946.727999999999
$ws.Cells.Item(31, 11).Value2 = 1726.7
$ws.Cells.Item(31, 12).Value2 = 9946.727999999999
$ws.Cells.Item(31, 13).Value2 = -1431.7
$ws.Cells.Item(31, 14).Value2 = -10536.728

$ws.Cells.Item(34, 8).Value2 = 6032.4287
$ws.Cells.Item(34, 9).Value2 = 1726.7
$ws.Cells.Item(34, 10).Value2 = 9946.727999999999
$ws.Cells.Item(34, 11).Value2 = 1726.7
$ws.Cells.Item(34, 12).Value2 = 9946.727999999999
$ws.Cells.Item(34, 13).Value2 = -1524.7
$ws.Cells.Item(34, 14).Value2 = -10350.728

$ws.Cells.Item(58, 8).Value2 = 2830.3572
$ws.Cells.Item(58, 9).Value2 = 1113.1
$ws.Cells.Item(58, 11).Value2 = 1113.1
$ws.Cells.Item(58, 13).Value2 = -910.0999999999999

$ws.Cells.Item(105, 8).Value2 = 2666.875
$ws.Cells.Item(105, 9).Value2 = 2619.2856
$ws.Cells.Item(105, 10).Value2 = 3000
$ws.Cells.Item(105, 11).Value2 = 2619.2856
$ws.Cells.Item(105, 12).Value2 = 3000
$ws.Cells.Item(105, 13).Value2 = -872.2856000000002
$ws.Cells.Item(105, 14).Value2 = -6494

$ws.Cells.Item(113, 8).Value2 = 952
$ws.Cells.Item(113, 9).Value2 = 502.5
$ws.Cells.Item(113, 11).Value2 = 502.5
$ws.Cells.Item(113, 13).Value2 = 1667.5

$ws.Cells.Item(122, 8).Value2 = 1556
$ws.Cells.Item(122, 9).Value2 = 1556
$ws.Cells.Item(122, 11).Value2 = 4668
$ws.Cells.Item(122, 13).Value2 = -2218

$ws.Cells.Item(132, 8).Value2 = 1757.0625
$ws.Cells.Item(132, 9).Value2 = 1757.0625
$ws.Cells.Item(132, 11).Value2 = 5271.1875
$ws.Cells.Item(132, 13).Value2 = -2741.1875

$ws.Cells.Item(136, 8).Value2 = 2830.3572
$ws.Cells.Item(136, 9).Value2 = 1113.1
$ws.Cells.Item(136, 11).Value2 = 3339.3
$ws.Cells.Item(136, 13).Value2 = -789.2999999999997

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value2 = 1141.7778
$ws.Cells.Item(5, 9).Value2 = 894.3333
$ws.Cells.Item(5, 10).Value2 = 1636.6666
$ws.Cells.Item(5, 11).Value2 = 2682.9999
$ws.Cells.Item(5, 12).Value2 = 4909.9998
$ws.Cells.Item(5, 13).Value2 = -2570.9999
$ws.Cells.Item(5, 14).Value2 = -5133.9998

$ws.Cells.Item(33, 8).Value2 = 288.375
$ws.Cells.Item(33, 9).Value2 = 355
$ws.Cells.Item(33, 10).Value2 = 88.5
$ws.Cells.Item(33, 11).Value2 = 2130
$ws.Cells.Item(33, 12).Value2 = 531
$ws.Cells.Item(33, 13).Value2 = -1847
$ws.Cells.Item(33, 14).Value2 = -1097

$ws.Cells.Item(38, 8).Value2 = 456.43478
$ws.Cells.Item(38, 9).Value2 = 441.41177
$ws.Cells.Item(38, 11).Value2 = 1324.23531
$ws.Cells.Item(38, 13).Value2 = -977.23531

$ws.Cells.Item(135, 8).Value2 = 1141.7778
$ws.Cells.Item(135, 9).Value2 = 894.3333
$ws.Cells.Item(135, 10).Value2 = 1636.6666
$ws.Cells.Item(135, 11).Value2 = 8048.9997
$ws.Cells.Item(135, 12).Value2 = 14729.9994
$ws.Cells.Item(135, 13).Value2 = -5513.9997
$ws.Cells.Item(135, 14).Value2 = -19799.9994

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value2 = 1813.2
$ws.Cells.Item(102, 9).Value2 = 1459.6666
$ws.Cells.Item(102, 10).Value2 = 4995
$ws.Cells.Item(102, 11).Value2 = 1459.6666
$ws.Cells.Item(102, 12).Value2 = 4995
$ws.Cells.Item(102, 13).Value2 = 162.3334
$ws.Cells.Item(102, 14).Value2 = -8239

$ws.Cells.Item(113, 8).Value2 = 3489.3333
$ws.Cells.Item(113, 9).Value2 = 1558.2
$ws.Cells.Item(113, 11).Value2 = 1558.2
$ws.Cells.Item(113, 13).Value2 = 611.8

$ws.Cells.Item(122, 8).Value2 = 2608.2
$ws.Cells.Item(122, 9).Value2 = 2608.2
$ws.Cells.Item(122, 11).Value2 = 7824.599999999999
$ws.Cells.Item(122, 13).Value2 = -5374.599999999999

$ws.Cells.Item(132, 8).Value2 = 102268
$ws.Cells.Item(132, 9).Value2 = 127022.5
$ws.Cells.Item(132, 10).Value2 = 3250
$ws.Cells.Item(132, 11).Value2 = 381067.5
$ws.Cells.Item(132, 12).Value2 = 9750
$ws.Cells.Item(132, 13).Value2 = -378537.5
$ws.Cells.Item(132, 14).Value2 = -14810

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value2 = 7898.7
$ws.Cells.Item(7, 9).Value2 = 7548.5
$ws.Cells.Item(7, 10).Value2 = 9299.5
$ws.Cells.Item(7, 11).Value2 = 7548.5
$ws.Cells.Item(7, 12).Value2 = 9299.5
$ws.Cells.Item(7, 13).Value2 = -7436.5
$ws.Cells.Item(7, 14).Value2 = -9523.5

$ws.Cells.Item(100, 8).Value2 = 8666.556
$ws.Cells.Item(100, 9).Value2 = 4000
$ws.Cells.Item(100, 10).Value2 = 9999.857
$ws.Cells.Item(100, 11).Value2 = 4000
$ws.Cells.Item(100, 12).Value2 = 9999.857
$ws.Cells.Item(100, 13).Value2 = -3459
$ws.Cells.Item(100, 14).Value2 = -11081.857

$ws.Cells.Item(122, 8).Value2 = 2749.25
$ws.Cells.Item(122, 9).Value2 = 2749.25
$ws.Cells.Item(122, 10).Value2 = 0
$ws.Cells.Item(122, 11).Value2 = 8247.75
$ws.Cells.Item(122, 12).Value2 = 0
$ws.Cells.Item(122, 13).ClearContents()
$ws.Cells.Item(122, 14).Value2 = -5797.75

$ws.Cells.Item(126, 8).Value2 = 7898.7
$ws.Cells.Item(126, 9).Value2 = 7548.5
$ws.Cells.Item(126, 10).Value2 = 9299.5
$ws.Cells.Item(126, 11).Value2 = 22645.5
$ws.Cells.Item(126, 12).Value2 = 27898.5
$ws.Cells.Item(126, 13).Value2 = -20175.5
$ws.Cells.Item(126, 14).Value2 = -32838.5

$ws.Cells.Item(132, 8).Value2 = 3348.75
$ws.Cells.Item(132, 9).Value2 = 3250
$ws.Cells.Item(132, 10).Value2 = 3447.5
$ws.Cells.Item(132, 11).Value2 = 9750
$ws.Cells.Item(132, 12).Value2 = 10342.5
$ws.Cells.Item(132, 13).Value2 = -7220
$ws.Cells.Item(132, 14).Value2 = -15402.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(59, 8).Value2 = 18000
$ws.Cells.Item(59, 10).Value2 = 18000
$ws.Cells.Item(59, 12).Value2 = 18000
$ws.Cells.Item(59, 14).Value2 = -19476

$ws.Cells.Item(122, 8).Value2 = 1447.7
$ws.Cells.Item(122, 9).Value2 = 1274.5
$ws.Cells.Item(122, 10).Value2 = 1563.1666
$ws.Cells.Item(122, 11).Value2 = 3823.5
$ws.Cells.Item(122, 12).Value2 = 4689.4998
$ws.Cells.Item(122, 13).Value2 = -1373.5
$ws.Cells.Item(122, 14).Value2 = -9589.4998

$ws.Cells.Item(126, 8).Value2 = 5790.0938
$ws.Cells.Item(126, 9).Value2 = 4722.8237
$ws.Cells.Item(126, 11).Value2 = 14168.4711
$ws.Cells.Item(126, 13).Value2 = -11698.4711

$ws.Cells.Item(136, 8).Value2 = 2642.24
$ws.Cells.Item(136, 9).Value2 = 2193.8635
$ws.Cells.Item(136, 11).Value2 = 6581.5905
$ws.Cells.Item(136, 13).Value2 = -4031.5905
